$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 50: 2024-04-19, 1 hour logged, running total continues from C49 ---
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)   # xlPasteFormats - reuse the date style from A49
$ws.Range("A50").Value = 45401
$ws.Range("B50").Value = 1

# --- Row 51: 2024-04-21, 1.5 hours logged, running total continues from C50 ---
$ws.Range("A49").Copy()
$ws.Range("A51").PasteSpecial(-4122)   # xlPasteFormats - reuse the date style from A49
$ws.Range("A51").Value = 45403
$ws.Range("B51").Value = 1.5

# Running-total column: fill C50:C51 in one shot so the new rows share one
# relative formula group (mirrors the existing C-column "running total"
# pattern already used throughout the sheet).
$ws.Range("C50:C51").Formula = "=C49+B50"

$excel.CutCopyMode = 0

$ws.Range("C51").Select()
